# Apply the "Working screenshot and link count" edit:
#  - Change the shared text value "Test" (used by the screenshot column, AH2:AH9)
#    to "imageUrl"
#  - Populate the "# of links" counts in column U (rows 2,3,4,5,6,8) which were
#    previously 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the screenshot column values (AH2:AH9) from "Test" to "imageUrl".
# All of these cells currently share the same text, so updating each one
# keeps them all in sync and (de-duplicated) results in the same shared
# string change shown in the diff.
foreach ($row in 2..9) {
    $cell = $ws.Cells.Item($row, 34)  # column AH = 34
    if ($cell.Value2 -eq "Test") {
        $cell.Value = "imageUrl"
    }
}

# Update link counts in column U (21)
$ws.Cells.Item(2, 21).Value = 28
$ws.Cells.Item(3, 21).Value = 118
$ws.Cells.Item(4, 21).Value = 3
$ws.Cells.Item(5, 21).Value = 2
$ws.Cells.Item(6, 21).Value = 75
$ws.Cells.Item(8, 21).Value = 3
